$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.342.04'
$ws.Range("E2").Value = '  +0.67%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.933.72'
$ws.Range("E3").Value = '  +1.21%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.74'
$ws.Range("E5").Value = '  +2.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7181'
$ws.Range("E6").Value = '  +1.32%  '
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3280'
$ws.Range("E8").Value = '  +1.70%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '27.76'
$ws.Range("E9").Value = '  +6.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07280'
$ws.Range("E10").Value = '  +6.69%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8054'
$ws.Range("E11").Value = '  +2.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08086'
$ws.Range("E12").Value = '  +2.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.934.92'
$ws.Range("E13").Value = '  +1.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.435'
$ws.Range("E14").Value = '  +1.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.64'
$ws.Range("E15").Value = '  +1.08%  '
$ws.Range("E16").Value = '  +4.94%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.336.43'
$ws.Range("E17").Value = '  +0.63%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008232'
$ws.Range("E18").Value = '  +5.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '253.41'
$ws.Range("E19").Value = '  -1.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.811'
$ws.Range("E20").Value = '  +0.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.188.09'
$ws.Range("E21").Value = '  +1.14%  '
$ws.Range("E22").Value = '  +0.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9959'
$ws.Range("E23").Value = '  -0.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.948'
$ws.Range("E24").Value = '  +2.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.727'
$ws.Range("E25").Value = '  +1.84%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.36'
$ws.Range("E26").Value = '  +4.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.343'
$ws.Range("E27").Value = '  +6.89%  '
$ws.Range("E28").Value = '  +3.26%  '
$ws.Range("E29").Value = '  -1.46%  '
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("E31").Value = '  +0.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.444'
$ws.Range("E32").Value = '  +1.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.207'
$ws.Range("E33").Value = '  +1.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05231'
$ws.Range("E34").Value = '  +4.73%  '
$ws.Range("E35").Value = '  +7.60%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7500'
$ws.Range("E36").Value = '  +1.74%  '
$ws.Range("E37").Value = '  +1.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01965'
$ws.Range("E38").Value = '  +2.27%  '
$ws.Range("E39").Value = '  +0.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '79.22'
$ws.Range("E40").Value = '  +0.24%  '
$ws.Range("E41").Value = '  +0.68%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4554'
$ws.Range("E42").Value = '  +3.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.034'
$ws.Range("E43").Value = '  +1.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8428'
$ws.Range("E44").Value = '  +1.53%  '
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.96'
$ws.Range("E46").Value = '  +0.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.812'
$ws.Range("E47").Value = '  +2.31%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.461'
$ws.Range("E48").Value = '  +3.67%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.79'
$ws.Range("E49").Value = '  +2.81%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4191'
$ws.Range("E50").Value = '  +3.74%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06054'
$ws.Range("E51").Value = '  +2.62%  '
